$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4664
$wsExhibit.Range("F3").Value = 1852
$wsExhibit.Range("F4").Value = 140
$wsExhibit.Range("F6").Value = 3127
$wsExhibit.Range("F8").Value = 591
$wsExhibit.Range("F9").Value = 272
$wsExhibit.Range("F10").Value = 635
$wsExhibit.Range("F11").Value = 542
$wsExhibit.Range("F12").Value = 535
$wsExhibit.Range("F13").Value = 388
$wsExhibit.Range("F15").Value = 1783
$wsExhibit.Range("F16").Value = 1348
$wsExhibit.Range("F17").Value = 126
$wsExhibit.Range("F18").Value = 1622
$wsExhibit.Range("F21").Value = 611
$wsExhibit.Range("F22").Value = 11
$wsExhibit.Range("F24").Value = 538
$wsExhibit.Range("F26").Value = 51
$wsExhibit.Range("F27").Value = 105
$wsExhibit.Range("F28").Value = 5
$wsExhibit.Range("F30").Value = 30
$wsExhibit.Range("F31").Value = 85
$wsExhibit.Range("F32").Value = 3877
$wsExhibit.Range("F33").Value = 4
$wsExhibit.Range("F34").Value = 767
$wsExhibit.Range("F36").Value = 936
$wsExhibit.Range("F38").Value = 1854

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 24
$wsShow.Range("F3").Value = 46

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4664
$wsAll.Range("F3").Value = 1852
$wsAll.Range("F4").Value = 140
$wsAll.Range("F6").Value = 3127
$wsAll.Range("F8").Value = 591
$wsAll.Range("F9").Value = 272
$wsAll.Range("F10").Value = 635
$wsAll.Range("F11").Value = 542
$wsAll.Range("F12").Value = 535
$wsAll.Range("F13").Value = 24
$wsAll.Range("F14").Value = 388
$wsAll.Range("F16").Value = 1783
$wsAll.Range("F17").Value = 1348
$wsAll.Range("F18").Value = 126
$wsAll.Range("F19").Value = 1622
$wsAll.Range("F22").Value = 611
$wsAll.Range("F23").Value = 11
$wsAll.Range("F25").Value = 538
$wsAll.Range("F27").Value = 51
$wsAll.Range("F28").Value = 105
$wsAll.Range("F29").Value = 5
$wsAll.Range("F31").Value = 30
$wsAll.Range("F32").Value = 85
$wsAll.Range("F33").Value = 3877
$wsAll.Range("F34").Value = 46
$wsAll.Range("F35").Value = 4
$wsAll.Range("F36").Value = 767
$wsAll.Range("F38").Value = 936
$wsAll.Range("F40").Value = 1854
